$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ordering: descending by value (row 12 "Malay-Indonesian" happens to stay in place)
$data = @(
    @("English", 27.41831963124954),
    @("Chinese", 8.291896622367046),
    @("Spanish", 8.018159801131798),
    @("Japanese", 7.019182309015775),
    @("German", 5.922189456597759),
    @("Arabic", 4.751012953768244),
    @("Portuguese", 3.697112823677367),
    @("French", 3.54707341997896),
    @("Italian", 3.324074735614117),
    @("Russian", 3.013375528608906),
    @("Malay-Indonesian", 2.531275167084562),
    @("Dutch", 1.663385872093),
    @("Korean", 1.501321332614216),
    @("Persian", 1.396834000522893),
    @("Turkish", 1.358933709332967),
    @("Thai", 0.9122042651057995),
    @("Polish", 0.8975563892822706),
    @("Urdu", 0.8109529189091431),
    @("Swedish", 0.5171063020517616),
    @("Bengali", 0.4287479724426412)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $row++
}

# Remove the now-unused trailing rows (previously held Uzbek and Vietnamese)
$ws.Range("A22:B23").Delete()
